# The revision moves all of this document's paragraph content elsewhere
# ("Chuyen gop y cua Liem vao phan chung"), leaving this file with a
# single empty paragraph. Clear the whole story down to one blank
# paragraph, same as selecting the entire body and pressing Delete.
$d = $word.ActiveDocument

$guard = 0
while (($d.Content.End - $d.Content.Start) -gt 1 -and $guard -lt 10) {
    $d.Content.Delete()
    $guard = $guard + 1
}

# Word also quietly marks the built-in "Normal Table" style as a quick
# style (<w:qFormat/>) during this kind of save.
$tableNormal = $d.Styles("Normal Table")
$tableNormal.QuickStyle = $true
